$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 88, pushing existing rows 88-107 down to 89-108.
$ws.Rows.Item(88).Insert()

# Populate the newly inserted row 88 with the new weekly record
# (Mercado Mayorista Lo Valledor de Santiago - Mora).
$ws.Range("A88").Value = 6
$ws.Range("B88").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C88").Value = "Metropolitana"
$ws.Range("D88").Value = 44932
$ws.Range("E88").Value = 13
$ws.Range("F88").Value = "Fruta"
$ws.Range("G88").Value = 100101
$ws.Range("H88").Value = "Berries"
$ws.Range("I88").Value = 100101008
$ws.Range("J88").Value = "Mora"
$ws.Range("K88").Value = "Sin especificar"
$ws.Range("L88").Value = "Primera"
$ws.Range("M88").Value = 300
$ws.Range("N88").Value = 4000
$ws.Range("O88").Value = 4000
$ws.Range("P88").Value = 4000
$ws.Range("Q88").Value = "$/bandeja 2 kilos"
$ws.Range("R88").Value = "Provincia de Curicó"
$ws.Range("S88").Value = 2000
$ws.Range("T88").Value = 2
